$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F23").Value = "trig mode"
$ws.Range("I23").ClearContents()

$ws.Range("E26").Value = 26
$ws.Range("F26").Value = "frequency"
$ws.Range("G26").Value = "dds_compiler"

$ws.Range("E27").Value = 25
$ws.Range("F27").Value = "fine_gain"
$ws.Range("G27").Value = "fine_gain"

$ws.Range("E28").Value = 2
$ws.Range("F28").Value = "coarse_gain"
$ws.Range("G28").Value = "coarse_gain_and_limiter"

$ws.Range("I26:I28").Select()
